# Add a new "2022-Q1" sheet (with per-fund holding detail) positioned
# between the existing "2021-Q4" sheet and the "总计" (summary) sheet, and
# add a corresponding "2022-Q1" row to the "总计" sheet.
#
# NOTE on ordering: PasteSpecial() relies on a pending Copy() clipboard, and
# *any* ClearFormats() call (even on an unrelated cell/sheet) silently drops
# that pending clipboard. So every PasteSpecial() in this script runs first,
# and every ClearFormats() call is deferred to the very end.

$wb = $excel.ActiveWorkbook

# --- locate the existing sheets -------------------------------------------------
$q4Sheet = $wb.Worksheets.Item("2021-Q4")

# style used for header cells / the leading index column - clone it from the
# existing "2021-Q4" sheet (its B1 cell) rather than hard-coding a style index.
$q4Sheet.Range("B1").Copy()

# --- insert the new "2022-Q1" sheet right after "2021-Q4" ----------------------
$q1Sheet = $wb.Worksheets.Add($null, $q4Sheet)
$q1Sheet.Name = "2022-Q1"

# --- header row: stamp the style on B1:H1 ---------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q1Sheet.Cells.Item(1, $i + 2).PasteSpecial(-4122)
}

# --- data rows: index, code, name, scale, total stock position, position pct, held value(亿元), position rank
$rows = @(
    @(0, "009983", "永赢港股通品质生活慧选混合",                         "9.75",  "87.97", "6.67", "0.6503", 1),
    @(1, "013991", "中欧港股通精选一年持有混合A",                        "12.87", "94.50", "4.33", "0.5573", 7),
    @(2, "011315", "永赢港股通优质成长一年持有期混合型证券投资基金",        "4.07",  "90.36", "6.27", "0.2552", 2),
    @(3, "013992", "中欧港股通精选一年持有混合C",                        "5.32",  "94.50", "4.33", "0.2304", 7),
    @(4, "008480", "永赢股息优选混合A",                                "3.42",  "89.21", "6.66", "0.2278", 3),
    @(5, "011203", "永赢惠添益混合A",                                  "5.60",  "93.07", "4.06", "0.2274", 8),
    @(6, "011071", "鹏华安悦一年持有期混合A",                           "9.16",  "21.81", "0.57", "0.0522", 9),
    @(7, "009140", "永赢竞争力精选混合",                                "1.51",  "88.85", "3.08", "0.0465", 8),
    @(8, "003413", "华泰柏瑞新经济沪港深灵活配置混合",                    "0.54",  "92.57", "5.76", "0.0311", 7),
    @(9, "008481", "永赢股息优选混合C",                                "0.19",  "89.21", "6.66", "0.0127", 3),
    @(10, "011204", "永赢惠添益混合C",                                 "0.29",  "93.07", "4.06", "0.0118", 8),
    @(11, "005493", "鑫元价值精选灵活配置混合A",                         "0.07",  "21.81", "0.57", "0.0004", 9)
)

# stamp the style on every leading index cell (A2:A13)
$r = 2
foreach ($row in $rows) {
    $q1Sheet.Cells.Item($r, 1).PasteSpecial(-4122)
    $r++
}

# --- "总计" sheet: insert a new data row for 2022-Q1 above the 2021-Q4 row -----
# (re-fetch by name now that sheet positions have shifted from inserting "2022-Q1")
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Cells.Item(2, 1).PasteSpecial(-4122)

# === all PasteSpecial calls are done from this point on - values + cleanup only ===

# --- header row values -----------------------------------------------------------
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q1Sheet.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# --- data row values ---------------------------------------------------------------
$r = 2
foreach ($row in $rows) {
    $q1Sheet.Cells.Item($r, 1).Value = $row[0]

    # Leading "'" forces these numeric-looking values to stay text (matching
    # the source data, which stores them as plain strings).
    $q1Sheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $q1Sheet.Cells.Item($r, 3).Value = $row[2]
    $q1Sheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $q1Sheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $q1Sheet.Cells.Item($r, 6).Value = "'" + $row[5]
    $q1Sheet.Cells.Item($r, 7).Value = "'" + $row[6]
    $q1Sheet.Cells.Item($r, 8).Value = $row[7]

    $r++
}

# --- "总计" sheet new row values ---------------------------------------------------
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 12
$totalSheet.Cells.Item(2, 4).Value = 2.3

# the previously-first data row (2021-Q4) shifted down to row 3 by the
# Insert() above, carrying its existing style along; just fix up its index.
$totalSheet.Cells.Item(3, 1).Value = 1

# --- cleanup: drop the quote-prefix / copied-format styling these edits picked up,
# so the plain data cells stay unstyled like the rest of the sheet (ClearFormats
# must come last - see note at top of file). ---------------------------------------
$q1Sheet.Range("B2:G13").ClearFormats()
$totalSheet.Range("B2:D2").ClearFormats()
